# feat: payment order validation and comments added

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Maximize the workbook window (matches xWindow=-120 yWindow=-120 windowWidth=38640 windowHeight=21240)
$excel.WindowState = -4137   # xlMaximized
$excel.Left = -120
$excel.Top = -120
$excel.Width = 38640
$excel.Height = 21240

# Update the username value cell (B1) text
$ws.Range("B1").Value = "jul3084dz1"

# Move the selection to D3
$ws.Range("D3").Select()
